$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.426.38'
$ws.Range("E2").Value = '  +2.30%  '
$ws.Range("D3").Value = '3.049.95'
$ws.Range("E3").Value = '  +2.46%  '
$ws.Range("D5").Value = '''516.70'
$ws.Range("E5").Value = '  +2.55%  '
$ws.Range("D6").Value = '''141.11'
$ws.Range("E6").Value = '  +3.02%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +1.60%  '
$ws.Range("D9").Value = '''7.24'
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").Value = '''0.373'
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("D12").Value = '3.580.54'
$ws.Range("E12").Value = '  +2.72%  '
$ws.Range("E13").Value = '  +3.22%  '
$ws.Range("D14").Value = '''25.47'
$ws.Range("E14").Value = '  -1.50%  '
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '57.484.01'
$ws.Range("E16").Value = '  +2.42%  '
$ws.Range("D17").Value = '3.055.62'
$ws.Range("E17").Value = '  +2.76%  '
$ws.Range("D18").Value = '''6.05'
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("D19").Value = '''12.76'
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").Value = '''8.08'
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("D21").Value = '''329.11'
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").Value = '''1.00'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '''0.494'
$ws.Range("E23").Value = '  +0.41%  '
$ws.Range("D24").Value = '''65.66'
$ws.Range("E24").Value = '  +1.95%  '
$ws.Range("E25").Value = '  +3.84%  '
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("D27").Value = '0.0₃0892'
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").Value = '''6.28'
$ws.Range("E28").Value = '  -0.98%  '
$ws.Range("D29").Value = '''7.11'
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("E30").Value = '  +2.59%  '
$ws.Range("D31").Value = '''20.60'
$ws.Range("E31").Value = '  +2.60%  '
$ws.Range("E32").Value = '  +2.20%  '
$ws.Range("D33").Value = '''154.39'
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("D34").Value = '''27.24'
$ws.Range("E34").Value = '  +6.13%  '
$ws.Range("D35").Value = '''4.46'
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("D36").Value = '''5.90'
$ws.Range("E36").Value = '  +2.60%  '
$ws.Range("D37").Value = '''1.26'
$ws.Range("E37").Value = '  +2.61%  '
$ws.Range("E38").Value = '  +1.74%  '
$ws.Range("D39").Value = '3.091.32'
$ws.Range("E39").Value = '  +2.69%  '
$ws.Range("D40").Value = '''3.88'
$ws.Range("E40").Value = '  +3.04%  '
$ws.Range("D41").Value = '''36.57'
$ws.Range("E41").Value = '  -0.74%  '
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").Value = '''0.653'
$ws.Range("E43").Value = '  +0.54%  '
$ws.Range("D44").Value = '2.247.75'
$ws.Range("E44").Value = '  +3.62%  '
$ws.Range("D45").Value = '''0.0256'
$ws.Range("E45").Value = '  +8.77%  '
$ws.Range("D46").Value = '''20.60'
$ws.Range("E46").Value = '  +6.18%  '
$ws.Range("D47").Value = '''1.35'
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("D48").Value = '''5.84'
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("D49").Value = '''0.914'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").Value = '''260.31'
$ws.Range("E50").Value = '  +15.38%  '
$ws.Range("D51").Value = '''0.710'
$ws.Range("E51").Value = '  +6.17%  '
